$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddEmployee")

# Update existing employee rows (replace Jane/Mark/James with new names)
$ws.Range("A3").Value = "Donald"
$ws.Range("C3").Value = "Trump"

$ws.Range("A4").Value = "Katie"
$ws.Range("C4").Value = "Ball"

$ws.Range("A5").Value = "Mohammed"
$ws.Range("C5").Value = "Salah"

# Add new Employee ID column
$ws.Range("D1").Value = "Employee ID"
$ws.Range("D2").Value = 55555555
$ws.Range("D3").Value = 66666666
$ws.Range("D4").Value = 3333333333
$ws.Range("D5").Value = 4444444

# Resize the new column to fit its content
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(4).ColumnWidth = 19.2

# Update the selected cell
$ws.Range("C11").Select()
